$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.024331898413565
$ws.Cells.Item(2, 4).Value = 1.029360826585827
$ws.Cells.Item(2, 5).Value = 1.035116203422241
$ws.Cells.Item(2, 6).Value = 1.04688578985722
$ws.Cells.Item(2, 9).Value = 1.031188489228352
$ws.Cells.Item(2, 10).Value = 1.029507107445881
$ws.Cells.Item(2, 11).Value = 1.032174962096279
$ws.Cells.Item(2, 12).Value = 1.037913731260087
$ws.Cells.Item(2, 13).Value = 1.04964997237128
$ws.Cells.Item(2, 14).Value = 1.013872143386366

$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.025188678055597
$ws.Cells.Item(3, 4).Value = 1.02998332774001
$ws.Cells.Item(3, 5).Value = 1.035906842015795
$ws.Cells.Item(3, 6).Value = 1.047828785796928
$ws.Cells.Item(3, 9).Value = 1.031309603444482
$ws.Cells.Item(3, 10).Value = 1.030003227856852
$ws.Cells.Item(3, 11).Value = 1.032606092400948
$ws.Cells.Item(3, 12).Value = 1.038513747836676
$ws.Cells.Item(3, 13).Value = 1.050404360396997
$ws.Cells.Item(3, 14).Value = 1.014038055163836

$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.025743692058625
$ws.Cells.Item(4, 4).Value = 1.030386417640582
$ws.Cells.Item(4, 5).Value = 1.036419388778886
$ws.Cells.Item(4, 6).Value = 1.048440091147144
$ws.Cells.Item(4, 9).Value = 1.031386769942654
$ws.Cells.Item(4, 10).Value = 1.030324236166947
$ws.Cells.Item(4, 11).Value = 1.032884661899132
$ws.Cells.Item(4, 12).Value = 1.03890229334802
$ws.Cells.Item(4, 13).Value = 1.050892997920137
$ws.Cells.Item(4, 14).Value = 1.014145357910424

$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.02597716696196
$ws.Cells.Item(5, 4).Value = 1.030555944325201
$ws.Cells.Item(5, 5).Value = 1.036635089338931
$ws.Cells.Item(5, 6).Value = 1.04869735121854
$ws.Cells.Item(5, 9).Value = 1.031418922206441
$ws.Cells.Item(5, 10).Value = 1.030459183327275
$ws.Cells.Item(5, 11).Value = 1.033001675337621
$ws.Cells.Item(5, 12).Value = 1.039065707069026
$ws.Cells.Item(5, 13).Value = 1.051098538796114
$ws.Cells.Item(5, 14).Value = 1.014190454784954

$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.026016377003714
$ws.Cells.Item(6, 4).Value = 1.030584412520174
$ws.Cells.Item(6, 5).Value = 1.036671319609468
$ws.Cells.Item(6, 6).Value = 1.048740561939935
$ws.Cells.Item(6, 9).Value = 1.031424303774426
$ws.Cells.Item(6, 10).Value = 1.03048184123693
$ws.Cells.Item(6, 11).Value = 1.033021316673439
$ws.Cells.Item(6, 12).Value = 1.039093148969667
$ws.Cells.Item(6, 13).Value = 1.051133056867795
$ws.Cells.Item(6, 14).Value = 1.014198025961537

$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.025746811185886
$ws.Cells.Item(7, 4).Value = 1.030388682600658
$ws.Cells.Item(7, 5).Value = 1.036422270094831
$ws.Cells.Item(7, 6).Value = 1.048443527620709
$ws.Cells.Item(7, 9).Value = 1.031387200697323
$ws.Cells.Item(7, 10).Value = 1.030326039357608
$ws.Cells.Item(7, 11).Value = 1.032886225821859
$ws.Cells.Item(7, 12).Value = 1.038904476619998
$ws.Cells.Item(7, 13).Value = 1.050895743907019
$ws.Cells.Item(7, 14).Value = 1.014145960549499

$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.024621321780036
$ws.Cells.Item(8, 4).Value = 1.029571142665966
$ws.Cells.Item(8, 5).Value = 1.035383205558066
$ws.Cells.Item(8, 6).Value = 1.047204245971983
$ws.Cells.Item(8, 9).Value = 1.031229668938038
$ws.Cells.Item(8, 10).Value = 1.029674776000099
$ws.Cells.Item(8, 11).Value = 1.032320746795599
$ws.Cells.Item(8, 12).Value = 1.038116447844313
$ws.Cells.Item(8, 13).Value = 1.04990481727432
$ws.Cells.Item(8, 14).Value = 1.013928224808669

$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.022642881814152
$ws.Cells.Item(9, 4).Value = 1.028132831096187
$ws.Cells.Item(9, 5).Value = 1.033559597866402
$ws.Cells.Item(9, 6).Value = 1.045029159440425
$ws.Cells.Item(9, 9).Value = 1.030942896972778
$ws.Cells.Item(9, 10).Value = 1.028527106549658
$ws.Cells.Item(9, 11).Value = 1.031321287992494
$ws.Cells.Item(9, 12).Value = 1.036730157155052
$ws.Cells.Item(9, 13).Value = 1.048162561556429
$ws.Cells.Item(9, 14).Value = 1.013544158426654

$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.021327247899774
$ws.Cells.Item(10, 4).Value = 1.027175605734728
$ws.Cells.Item(10, 5).Value = 1.032348907215044
$ws.Cells.Item(10, 6).Value = 1.043585046947301
$ws.Cells.Item(10, 9).Value = 1.030745581778463
$ws.Cells.Item(10, 10).Value = 1.027762024314585
$ws.Cells.Item(10, 11).Value = 1.030653031970745
$ws.Cells.Item(10, 12).Value = 1.035807601634606
$ws.Cells.Item(10, 13).Value = 1.047003762155801
$ws.Cells.Item(10, 14).Value = 1.013287879803995

$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.020758371197264
$ws.Cells.Item(11, 4).Value = 1.026761529617242
$ws.Cells.Item(11, 5).Value = 1.031825883463704
$ws.Cells.Item(11, 6).Value = 1.042961162524129
$ws.Cells.Item(11, 9).Value = 1.030658695402484
$ws.Cells.Item(11, 10).Value = 1.027430757521211
$ws.Cells.Item(11, 11).Value = 1.030363223034945
$ws.Cells.Item(11, 12).Value = 1.035408530454174
$ws.Cells.Item(11, 13).Value = 1.046502650285169
$ws.Cells.Item(11, 14).Value = 1.013176858344531

$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.020547186937694
$ws.Cells.Item(12, 4).Value = 1.02660778649437
$ws.Cells.Item(12, 5).Value = 1.031631793166013
$ws.Cells.Item(12, 6).Value = 1.042729640061866
$ws.Cells.Item(12, 9).Value = 1.030626205087657
$ws.Cells.Item(12, 10).Value = 1.027307714316326
$ws.Cells.Item(12, 11).Value = 1.030255508802248
$ws.Cells.Item(12, 12).Value = 1.035260359324093
$ws.Cells.Item(12, 13).Value = 1.04631661514144
$ws.Cells.Item(12, 14).Value = 1.013135612799033

$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.020592481139486
$ws.Cells.Item(13, 4).Value = 1.02664076202854
$ws.Cells.Item(13, 5).Value = 1.031673417839086
$ws.Cells.Item(13, 6).Value = 1.042779292613323
$ws.Cells.Item(13, 9).Value = 1.030633184172962
$ws.Cells.Item(13, 10).Value = 1.027334107299771
$ws.Cells.Item(13, 11).Value = 1.030278616850235
$ws.Cells.Item(13, 12).Value = 1.035292139721681
$ws.Cells.Item(13, 13).Value = 1.046356515755503
$ws.Cells.Item(13, 14).Value = 1.013144460428627

$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.020740912145742
$ws.Cells.Item(14, 4).Value = 1.026748819863591
$ws.Cells.Item(14, 5).Value = 1.031809836135384
$ws.Cells.Item(14, 6).Value = 1.042942020388731
$ws.Cells.Item(14, 9).Value = 1.030656014162824
$ws.Cells.Item(14, 10).Value = 1.027420586642335
$ws.Cells.Item(14, 11).Value = 1.030354320687841
$ws.Cells.Item(14, 12).Value = 1.03539628131631
$ws.Cells.Item(14, 13).Value = 1.046487270510694
$ws.Cells.Item(14, 14).Value = 1.013173449119815

$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.020832381581341
$ws.Cells.Item(15, 4).Value = 1.026815406287473
$ws.Cells.Item(15, 5).Value = 1.031893912385059
$ws.Cells.Item(15, 6).Value = 1.043042311015983
$ws.Cells.Item(15, 9).Value = 1.030670051757959
$ws.Cells.Item(15, 10).Value = 1.027473869989557
$ws.Cells.Item(15, 11).Value = 1.030400955574537
$ws.Cells.Item(15, 12).Value = 1.035460454603147
$ws.Cells.Item(15, 13).Value = 1.046567846144376
$ws.Cells.Item(15, 14).Value = 1.013191309064279

$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.021365019383952
$ws.Cells.Item(16, 4).Value = 1.027203095366547
$ws.Cells.Item(16, 5).Value = 1.032383644289024
$ws.Cells.Item(16, 6).Value = 1.043626482293794
$ws.Cells.Item(16, 9).Value = 1.03075131769641
$ws.Cells.Item(16, 10).Value = 1.027784009900329
$ws.Cells.Item(16, 11).Value = 1.030672256269554
$ws.Cells.Item(16, 12).Value = 1.035834095248835
$ws.Cells.Item(16, 13).Value = 1.047037033253923
$ws.Cells.Item(16, 14).Value = 1.013295246893789

$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.021699344613909
$ws.Cells.Item(17, 4).Value = 1.027446393175342
$ws.Cells.Item(17, 5).Value = 1.032691166226032
$ws.Cells.Item(17, 6).Value = 1.043993300488813
$ws.Cells.Item(17, 9).Value = 1.030801906560836
$ws.Cells.Item(17, 10).Value = 1.027978558368328
$ws.Cells.Item(17, 11).Value = 1.030842316533762
$ws.Cells.Item(17, 12).Value = 1.036068578594762
$ws.Cells.Item(17, 13).Value = 1.047331518649419
$ws.Cells.Item(17, 14).Value = 1.013360430977567

$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.021894428095371
$ws.Cells.Item(18, 4).Value = 1.027588343993865
$ws.Cells.Item(18, 5).Value = 1.03287065551805
$ws.Cells.Item(18, 6).Value = 1.044207396787087
$ws.Cells.Item(18, 9).Value = 1.030831274474776
$ws.Cells.Item(18, 10).Value = 1.028092036937895
$ws.Cells.Item(18, 11).Value = 1.030941466337224
$ws.Cells.Item(18, 12).Value = 1.036205387405834
$ws.Cells.Item(18, 13).Value = 1.047503350179185
$ws.Cells.Item(18, 14).Value = 1.013398446787949

$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.02196095957634
$ws.Cells.Item(19, 4).Value = 1.027636752170626
$ws.Cells.Item(19, 5).Value = 1.032931876512215
$ws.Cells.Item(19, 6).Value = 1.044280421352097
$ws.Cells.Item(19, 9).Value = 1.030841264457031
$ws.Cells.Item(19, 10).Value = 1.02813073045621
$ws.Cells.Item(19, 11).Value = 1.030975266445863
$ws.Cells.Item(19, 12).Value = 1.036252042187968
$ws.Cells.Item(19, 13).Value = 1.047561950959054
$ws.Cells.Item(19, 14).Value = 1.01341140834143

$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.021663466681787
$ws.Cells.Item(20, 4).Value = 1.027420285542879
$ws.Cells.Item(20, 5).Value = 1.032658159923489
$ws.Cells.Item(20, 6).Value = 1.043953930137959
$ws.Cells.Item(20, 9).Value = 1.030796493301199
$ws.Cells.Item(20, 10).Value = 1.027957684970307
$ws.Cells.Item(20, 11).Value = 1.030824075164153
$ws.Cells.Item(20, 12).Value = 1.03604341674406
$ws.Cells.Item(20, 13).Value = 1.047299916607012
$ws.Cells.Item(20, 14).Value = 1.013353437852366

$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.020697199537752
$ws.Cells.Item(21, 4).Value = 1.026716997766061
$ws.Cells.Item(21, 5).Value = 1.031769659253776
$ws.Cells.Item(21, 6).Value = 1.042894095138057
$ws.Cells.Item(21, 9).Value = 1.030649297280715
$ws.Cells.Item(21, 10).Value = 1.027395120510472
$ws.Cells.Item(21, 11).Value = 1.030332029614293
$ws.Cells.Item(21, 12).Value = 1.035365612498881
$ws.Cells.Item(21, 13).Value = 1.046448763717906
$ws.Cells.Item(21, 14).Value = 1.013164912868828

$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.020090374904402
$ws.Cells.Item(22, 4).Value = 1.02627517926354
$ws.Cells.Item(22, 5).Value = 1.031212089105446
$ws.Cells.Item(22, 6).Value = 1.042228985347866
$ws.Cells.Item(22, 9).Value = 1.030555495012204
$ws.Cells.Item(22, 10).Value = 1.027041437826475
$ws.Cells.Item(22, 11).Value = 1.030022278140122
$ws.Cells.Item(22, 12).Value = 1.034939807574845
$ws.Cells.Item(22, 13).Value = 1.045914190140479
$ws.Cells.Item(22, 14).Value = 1.013046338107578

$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.020411996725068
$ws.Cells.Item(23, 4).Value = 1.026509360287313
$ws.Cells.Item(23, 5).Value = 1.031507566017817
$ws.Cells.Item(23, 6).Value = 1.042581453500925
$ws.Cells.Item(23, 9).Value = 1.030605340048076
$ws.Cells.Item(23, 10).Value = 1.027228929041786
$ws.Cells.Item(23, 11).Value = 1.030186519129473
$ws.Cells.Item(23, 12).Value = 1.035165500479164
$ws.Cells.Item(23, 13).Value = 1.046197522202899
$ws.Cells.Item(23, 14).Value = 1.013109200627541

$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.021679678136497
$ws.Cells.Item(24, 4).Value = 1.027432082336919
$ws.Cells.Item(24, 5).Value = 1.032673073689728
$ws.Cells.Item(24, 6).Value = 1.043971719480338
$ws.Cells.Item(24, 9).Value = 1.030798939751987
$ws.Cells.Item(24, 10).Value = 1.027967116755155
$ws.Cells.Item(24, 11).Value = 1.030832317788359
$ws.Cells.Item(24, 12).Value = 1.036054786182076
$ws.Cells.Item(24, 13).Value = 1.047314196014645
$ws.Cells.Item(24, 14).Value = 1.013356597760017

$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.023153776179218
$ws.Cells.Item(25, 4).Value = 1.028504386052312
$ws.Cells.Item(25, 5).Value = 1.03403016188811
$ws.Cells.Item(25, 6).Value = 1.045590431272697
$ws.Cells.Item(25, 9).Value = 1.031018118384643
$ws.Cells.Item(25, 10).Value = 1.028823806409782
$ws.Cells.Item(25, 11).Value = 1.031580021287761
$ws.Cells.Item(25, 12).Value = 1.037088263677518
$ws.Cells.Item(25, 13).Value = 1.048612506586321
$ws.Cells.Item(25, 14).Value = 1.013643492049668
